# 20190730 before done, 20190731 init
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Finish filling in row 26 (2019-07-25): mark the previously-blank
#     columns B, D, H as "X" (not done) ---
$ws.Range("B26").Value = "×"
$ws.Range("D26").Value = "×"
$ws.Range("H26").Value = "×"

# --- Row 27: 2019-07-26 ---
$ws.Range("A27").Value = 20190726
$ws.Range("C27").Value = "√"
$ws.Range("D27").Value = "×"
$ws.Range("E27").Value = "√"
$ws.Range("F27").Value = "√"
$ws.Range("G27").Value = "√"
$ws.Range("H27").Value = "×"
$ws.Range("I27").Value = "√"
$ws.Range("J27").Value = "×"
$ws.Range("K27").Value = "√"
$ws.Range("L27").Value = "√"
$ws.Range("M27").Value = "√"
$ws.Range("N27").Value = "√"
$ws.Range("O27").Value = "×"

# --- Row 28: 2019-07-27 ---
$ws.Range("A28").Value = 20190727
$ws.Range("B28").Value = "×"
$ws.Range("C28").Value = "×"
$ws.Range("D28").Value = "×"
$ws.Range("E28").Value = "×"
$ws.Range("F28").Value = "√"
$ws.Range("G28").Value = "√"
$ws.Range("H28").Value = "×"
$ws.Range("I28").Value = "√"
$ws.Range("J28").Value = "×"
$ws.Range("K28").Value = "√"
$ws.Range("L28").Value = "√"
$ws.Range("M28").Value = "√"
$ws.Range("N28").Value = "√"
$ws.Range("O28").Value = "×"

# --- Row 29: 2019-07-28 ---
$ws.Range("A29").Value = 20190728
$ws.Range("B29").Value = "×"
$ws.Range("C29").Value = "×"
$ws.Range("D29").Value = "×"
$ws.Range("E29").Value = "×"
$ws.Range("F29").Value = "√"
$ws.Range("G29").Value = "√"
$ws.Range("H29").Value = "×"
$ws.Range("I29").Value = "√"
$ws.Range("J29").Value = "×"
$ws.Range("K29").Value = "√"
$ws.Range("L29").Value = "√"
$ws.Range("M29").Value = "√"
$ws.Range("N29").Value = "√"
$ws.Range("O29").Value = "×"

# --- Row 30: 2019-07-29 ---
$ws.Range("A30").Value = 20190729
$ws.Range("B30").Value = "√"
$ws.Range("C30").Value = "√"
$ws.Range("D30").Value = "√"
$ws.Range("E30").Value = "√"
$ws.Range("F30").Value = "√"
$ws.Range("G30").Value = "√"
$ws.Range("H30").Value = "√"
$ws.Range("I30").Value = "√"
$ws.Range("J30").Value = "√"
$ws.Range("K30").Value = "√"
$ws.Range("L30").Value = "√"
$ws.Range("M30").Value = "√"
$ws.Range("N30").Value = "√"
$ws.Range("O30").Value = "√"

# --- Row 31: 2019-07-30 ---
$ws.Range("A31").Value = 20190730
$ws.Range("B31").Value = "√"
$ws.Range("C31").Value = "√"
$ws.Range("D31").Value = "×"
$ws.Range("E31").Value = "√"
$ws.Range("F31").Value = "√"
$ws.Range("G31").Value = "√"
$ws.Range("H31").Value = "×"
$ws.Range("I31").Value = "√"
$ws.Range("J31").Value = "√"
$ws.Range("K31").Value = "√"
$ws.Range("L31").Value = "√"
$ws.Range("M31").Value = "√"
$ws.Range("N31").Value = "√"
$ws.Range("O31").Value = "√"

# --- Row 32: 2019-07-31 (just initialized with the date) ---
$ws.Range("A32").Value = 20190731

# Move the active selection to B32, matching the author's last position.
[void]$ws.Range("B32").Select()
